$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 (pushes existing rows 13:25 down to 14:26,
# carrying formats/values with them; dimension grows to A1:R26).
$ws.Rows(13).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Vega Modelo de Temuco"
$ws.Range("C13").Value = "La Araucanía"
$ws.Range("D13").Value = 44467
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 100112010
$ws.Range("G13").Value = "Achicoria"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("N13").Value = "$/caja 16 unidades"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 438
$ws.Range("Q13").Value = 16
$ws.Range("R13").Value = "Hortaliza"
